# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the zh-cn and
# de-de report sheets, row 2, to reflect the regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 16:51:56"
$wsZhCn.Range("H2").Value = "2016-03-17 16:52:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 16:52:00"
$wsDeDe.Range("H2").Value = "2016-03-17 16:52:19"
